$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.063857049188834
$ws.Cells.Item(2, 4).Value = 1.062743240390809
$ws.Cells.Item(2, 5).Value = 1.0686520845984
$ws.Cells.Item(2, 6).Value = 1.078065579570892
$ws.Cells.Item(2, 9).Value = 1.05206115945295
$ws.Cells.Item(2, 10).Value = 1.068820272934943
$ws.Cells.Item(2, 11).Value = 1.065463760880149
$ws.Cells.Item(2, 12).Value = 1.071356694158912
$ws.Cells.Item(2, 13).Value = 1.080745230416667

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.065093264452941
$ws.Cells.Item(3, 4).Value = 1.063711229497852
$ws.Cells.Item(3, 5).Value = 1.06976363419444
$ws.Cells.Item(3, 6).Value = 1.079280221981994
$ws.Cells.Item(3, 9).Value = 1.052454004787267
$ws.Cells.Item(3, 10).Value = 1.069710114300827
$ws.Cells.Item(3, 11).Value = 1.066246240683858
$ws.Cells.Item(3, 12).Value = 1.072283528993062
$ws.Cells.Item(3, 13).Value = 1.081776717094959

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.065892923662332
$ws.Cells.Item(4, 4).Value = 1.064337270394392
$ws.Cells.Item(4, 5).Value = 1.070482917961386
$ws.Cells.Item(4, 6).Value = 1.08006639377219
$ws.Cells.Item(4, 9).Value = 1.052706863294401
$ws.Cells.Item(4, 10).Value = 1.070285102132902
$ws.Cells.Item(4, 11).Value = 1.066751623846446
$ws.Cells.Item(4, 12).Value = 1.072882696669001
$ws.Cells.Item(4, 13).Value = 1.082443799659464

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.066229041488042
$ws.Cells.Item(5, 4).Value = 1.064600384040477
$ws.Cells.Item(5, 5).Value = 1.070785314909473
$ws.Cells.Item(5, 6).Value = 1.080396953763108
$ws.Cells.Item(5, 9).Value = 1.052812845102038
$ws.Cells.Item(5, 10).Value = 1.070526636950646
$ws.Cells.Item(5, 11).Value = 1.066963864412703
$ws.Cells.Item(5, 12).Value = 1.073134454712385
$ws.Cells.Item(5, 13).Value = 1.082724156683314

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.066285473756094
$ws.Cells.Item(6, 4).Value = 1.064644557670262
$ws.Cells.Item(6, 5).Value = 1.070836089312904
$ws.Cells.Item(6, 6).Value = 1.080452459411557
$ws.Cells.Item(6, 9).Value = 1.052830621171062
$ws.Cells.Item(6, 10).Value = 1.070567180599297
$ws.Cells.Item(6, 11).Value = 1.066999487518011
$ws.Cells.Item(6, 12).Value = 1.073176718286176
$ws.Cells.Item(6, 13).Value = 1.082771224944406

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.065897415114028
$ws.Cells.Item(7, 4).Value = 1.064340786421203
$ws.Cells.Item(7, 5).Value = 1.070486958563169
$ws.Cells.Item(7, 6).Value = 1.080070810519459
$ws.Cells.Item(7, 9).Value = 1.052708280684803
$ws.Cells.Item(7, 10).Value = 1.070288330278263
$ws.Cells.Item(7, 11).Value = 1.066754460688379
$ws.Cells.Item(7, 12).Value = 1.072886061189702
$ws.Cells.Item(7, 13).Value = 1.082447546134101

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.064274887918895
$ws.Cells.Item(8, 4).Value = 1.06307044216866
$ws.Cells.Item(8, 5).Value = 1.069027731099452
$ws.Cells.Item(8, 6).Value = 1.078476029609323
$ws.Cells.Item(8, 9).Value = 1.052194200905632
$ws.Cells.Item(8, 10).Value = 1.069121165172785
$ws.Cells.Item(8, 11).Value = 1.065728397495772
$ws.Cells.Item(8, 12).Value = 1.071670038246082
$ws.Cells.Item(8, 13).Value = 1.0810939011316

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.061413743593205
$ws.Cells.Item(9, 4).Value = 1.060829495441937
$ws.Cells.Item(9, 5).Value = 1.066456612621051
$ws.Cells.Item(9, 6).Value = 1.075667418608601
$ws.Cells.Item(9, 9).Value = 1.051278047836727
$ws.Cells.Item(9, 10).Value = 1.067058296391802
$ws.Cells.Item(9, 11).Value = 1.063913148223365
$ws.Cells.Item(9, 12).Value = 1.069522936318202
$ws.Cells.Item(9, 13).Value = 1.078705803555414

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.059504807249836
$ws.Cells.Item(10, 4).Value = 1.059333830171814
$ws.Cells.Item(10, 5).Value = 1.064742605969384
$ws.Cells.Item(10, 6).Value = 1.073795991791328
$ws.Cells.Item(10, 9).Value = 1.050660330551296
$ws.Cells.Item(10, 10).Value = 1.065678819800986
$ws.Cells.Item(10, 11).Value = 1.062698080156298
$ws.Cells.Item(10, 12).Value = 1.068088558832393
$ws.Cells.Item(10, 13).Value = 1.077111770974687

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.058677829183443
$ws.Cells.Item(11, 4).Value = 1.058685772652102
$ws.Cells.Item(11, 5).Value = 1.064000420694308
$ws.Cells.Item(11, 6).Value = 1.072985854376029
$ws.Cells.Item(11, 9).Value = 1.050391195233523
$ws.Cells.Item(11, 10).Value = 1.065080469426528
$ws.Cells.Item(11, 11).Value = 1.06217076458821
$ws.Cells.Item(11, 12).Value = 1.067466734095874
$ws.Cells.Item(11, 13).Value = 1.076421051263112

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.058370590515273
$ws.Cells.Item(12, 4).Value = 1.058444990092801
$ws.Cells.Item(12, 5).Value = 1.063724736970102
$ws.Cells.Item(12, 6).Value = 1.072684961777557
$ws.Cells.Item(12, 9).Value = 1.050290976223825
$ws.Cells.Item(12, 10).Value = 1.064858059353597
$ws.Cells.Item(12, 11).Value = 1.061974716978239
$ws.Cells.Item(12, 12).Value = 1.067235649451026
$ws.Cells.Item(12, 13).Value = 1.076164411365345

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.058436497110218
$ws.Cells.Item(13, 4).Value = 1.058496641738697
$ws.Cells.Item(13, 5).Value = 1.063783872230966
$ws.Cells.Item(13, 6).Value = 1.072749503002633
$ws.Cells.Item(13, 9).Value = 1.050312484875556
$ws.Cells.Item(13, 10).Value = 1.06490577416006
$ws.Cells.Item(13, 11).Value = 1.062016777984339
$ws.Cells.Item(13, 12).Value = 1.067285222951919
$ws.Cells.Item(13, 13).Value = 1.07621946495498

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.058652434023604
$ws.Cells.Item(14, 4).Value = 1.058665870827347
$ws.Cells.Item(14, 5).Value = 1.063977632664099
$ws.Cells.Item(14, 6).Value = 1.072960981929704
$ws.Cells.Item(14, 9).Value = 1.050382916207071
$ws.Cells.Item(14, 10).Value = 1.065062088134635
$ws.Cells.Item(14, 11).Value = 1.062154562893819
$ws.Cells.Item(14, 12).Value = 1.067447634845819
$ws.Cells.Item(14, 13).Value = 1.076399838886475

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.058785471560915
$ws.Cells.Item(15, 4).Value = 1.058770129781686
$ws.Cells.Item(15, 5).Value = 1.064097014404407
$ws.Cells.Item(15, 6).Value = 1.073091284795619
$ws.Cells.Item(15, 9).Value = 1.050426278100382
$ws.Cells.Item(15, 10).Value = 1.065158377598268
$ws.Cells.Item(15, 11).Value = 1.062239432946161
$ws.Cells.Item(15, 12).Value = 1.067547687384048
$ws.Cells.Item(15, 13).Value = 1.076510963116101

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.059559682486749
$ws.Cells.Item(16, 4).Value = 1.059376830554194
$ws.Cells.Item(16, 5).Value = 1.064791862025169
$ws.Cells.Item(16, 6).Value = 1.073849762048636
$ws.Cells.Item(16, 9).Value = 1.050678157138526
$ws.Cells.Item(16, 10).Value = 1.065718508545366
$ws.Cells.Item(16, 11).Value = 1.062733051309649
$ws.Cells.Item(16, 12).Value = 1.068129811778961
$ws.Cells.Item(16, 13).Value = 1.077157601228066

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.060045216616109
$ws.Cells.Item(17, 4).Value = 1.059757283292209
$ws.Cells.Item(17, 5).Value = 1.065227717947015
$ws.Cells.Item(17, 6).Value = 1.074325587694898
$ws.Cells.Item(17, 9).Value = 1.050835709285124
$ws.Cells.Item(17, 10).Value = 1.066069587511151
$ws.Cells.Item(17, 11).Value = 1.063042367347255
$ws.Cells.Item(17, 12).Value = 1.068494766507408
$ws.Cells.Item(17, 13).Value = 1.077563086975599

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.060328382431119
$ws.Cells.Item(18, 4).Value = 1.059979154063713
$ws.Cells.Item(18, 5).Value = 1.065481944572083
$ws.Cells.Item(18, 6).Value = 1.074603148439399
$ws.Cells.Item(18, 9).Value = 1.050927446710176
$ws.Cells.Item(18, 10).Value = 1.066274266868012
$ws.Cells.Item(18, 11).Value = 1.063222672103732
$ws.Cells.Item(18, 12).Value = 1.068707568118528
$ws.Cells.Item(18, 13).Value = 1.077799552496662

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.060424928287103
$ws.Cells.Item(19, 4).Value = 1.060054799370613
$ws.Cells.Item(19, 5).Value = 1.065568629203464
$ws.Cells.Item(19, 6).Value = 1.0746977928796
$ws.Cells.Item(19, 9).Value = 1.050958699674368
$ws.Cells.Item(19, 10).Value = 1.066344040486407
$ws.Cells.Item(19, 11).Value = 1.063284132061961
$ws.Cells.Item(19, 12).Value = 1.068780116091314
$ws.Cells.Item(19, 13).Value = 1.0778801731321

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.059993127353213
$ws.Cells.Item(20, 4).Value = 1.059716468552456
$ws.Cells.Item(20, 5).Value = 1.065180954823459
$ws.Cells.Item(20, 6).Value = 1.074274534084269
$ws.Cells.Item(20, 9).Value = 1.050818821994782
$ws.Cells.Item(20, 10).Value = 1.066031930320374
$ws.Cells.Item(20, 11).Value = 1.063009192471852
$ws.Cells.Item(20, 12).Value = 1.068455617603699
$ws.Cells.Item(20, 13).Value = 1.077519587114633

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.058588847723192
$ws.Cells.Item(21, 4).Value = 1.058616038887782
$ws.Cells.Item(21, 5).Value = 1.063920575144264
$ws.Cells.Item(21, 6).Value = 1.072898705885368
$ws.Cells.Item(21, 9).Value = 1.050362182847926
$ws.Cells.Item(21, 10).Value = 1.065016061887774
$ws.Cells.Item(21, 11).Value = 1.062113993628167
$ws.Cells.Item(21, 12).Value = 1.06739981167505
$ws.Cells.Item(21, 13).Value = 1.076346725367154

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.057705560216798
$ws.Cells.Item(22, 4).Value = 1.057923777518498
$ws.Cells.Item(22, 5).Value = 1.063128105822981
$ws.Cells.Item(22, 6).Value = 1.072033832064472
$ws.Cells.Item(22, 9).Value = 1.050073627880074
$ws.Cells.Item(22, 10).Value = 1.064376441094569
$ws.Cells.Item(22, 11).Value = 1.061550109509574
$ws.Cells.Item(22, 12).Value = 1.066735340122136
$ws.Cells.Item(22, 13).Value = 1.075608860958362

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.05817384277179
$ws.Cells.Item(23, 4).Value = 1.058290794610076
$ws.Cells.Item(23, 5).Value = 1.063548211136577
$ws.Cells.Item(23, 6).Value = 1.072492303151141
$ws.Cells.Item(23, 9).Value = 1.050226733839834
$ws.Cells.Item(23, 10).Value = 1.06471560242604
$ws.Cells.Item(23, 11).Value = 1.061849133916374
$ws.Cells.Item(23, 12).Value = 1.067087650688795
$ws.Cells.Item(23, 13).Value = 1.076000059108402

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.060016664366464
$ws.Cells.Item(24, 4).Value = 1.05973491110289
$ws.Cells.Item(24, 5).Value = 1.065202085068531
$ws.Cells.Item(24, 6).Value = 1.074297602952558
$ws.Cells.Item(24, 9).Value = 1.050826453129025
$ws.Cells.Item(24, 10).Value = 1.06604894629009
$ws.Cells.Item(24, 11).Value = 1.063024183122517
$ws.Cells.Item(24, 12).Value = 1.068473307525438
$ws.Cells.Item(24, 13).Value = 1.077539242976066

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.062153673845778
$ws.Cells.Item(25, 4).Value = 1.061409129120434
$ws.Cells.Item(25, 5).Value = 1.06712128974018
$ws.Cells.Item(25, 6).Value = 1.076393332219230
$ws.Cells.Item(25, 9).Value = 1.051516117333636
$ws.Cells.Item(25, 10).Value = 1.06759233729853
$ws.Cells.Item(25, 11).Value = 1.064383292923894
$ws.Cells.Item(25, 12).Value = 1.070078532653069
$ws.Cells.Item(25, 13).Value = 1.079323524676658
